$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the "5000 process" results table (rows 46-82, columns A-H) while
# keeping cell formatting intact. This removes the "Tempo médio"/"Tempo
# total" section headers as well as all of the associated benchmark
# data/formulas underneath them.
$ws.Range("A46:H82").ClearContents()

# A new (narrower) column was inserted into the column-width list at AZ
# (column 52), matching the author's manual resize of that column.
$ws.Columns.Item(52).ColumnWidth = 10.5703125

# The author's selection/scroll position moved when finishing the edit.
$ws.Application.ActiveWindow.ScrollColumn = $ws.Range("AD1").Column
$ws.Range("AI6").Select()
